# DaySale_2025-08-07_00-00.xlsx — "Upload new version with timestamp"
#
# The original sheet has a single sold-item row (row 7), a totals row
# (old row 8, P8:Q8 merged) and a footer row (old row 9: timestamp /
# "1/1" / "developed by" merged blocks).
#
# The new version adds a second sold item ("POWER B COMPLEX I.M./I.V.
# 6 AMP") as a new row 8 — styled exactly like row 7 — which pushes the
# totals row down to row 9 (new total 38.32) and the footer row down to
# row 10 (with a refreshed timestamp, 9:30 AM instead of 9:22 AM).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a fresh blank row at 8: old row 8 (totals) -> row 9,
#    old row 9 (footer) -> row 10. Excel shifts the existing merged
#    ranges (P8:Q8 -> P9:Q9, A9:F9 -> A10:F10, etc.) automatically.
$ws.Rows(8).EntireRow.Insert()

# 2. Clone row 7's formatting (fonts/fills/borders/number formats) onto
#    the new row 8 so the new item row looks identical to the first one.
$ws.Range("A7:Q7").Copy()
$ws.Range("A8:Q8").PasteSpecial(-4122)  # xlPasteFormats
$ws.Application.CutCopyMode = $false

# 3. Re-create row 7's merge layout on row 8.
$ws.Range("A8:B8").Merge()
$ws.Range("C8:G8").Merge()
$ws.Range("H8:K8").Merge()
$ws.Range("L8:M8").Merge()
$ws.Range("N8:O8").Merge()

# 4. Row heights: new item row keeps the old totals row's height
#    (24.75), the totals row (now row 9) grows slightly to 25.5.
$ws.Rows(8).RowHeight = 24.75
$ws.Rows(9).RowHeight = 25.5

# 5. Fill in the new item's data (row 8).
$ws.Range("A8").Value = 2
$ws.Range("C8").Value = "POWER B COMPLEX I.M./I.V. 6 AMP"
$ws.Range("H8").Value = "0:3"

# L8/N8/P8/Q8 hold numeric-looking text (same as row 7), so force text
# storage via a temporary "@" number format, then restore the original
# numeric display format copied from row 7.
$ws.Range("L8").NumberFormat = "@"
$ws.Range("L8").Value = "1"
$ws.Range("L8").NumberFormat = "#,##0.##;""[""#,##0.##""]"";0"

$ws.Range("N8").Value = "48.00"

$ws.Range("P8").NumberFormat = "@"
$ws.Range("P8").Value = "-7.6800"
$ws.Range("P8").NumberFormat = "0.00"

$ws.Range("Q8").Value = "0:-1"

# 6. Update the footer timestamp (now row 10).
$ws.Range("A10").Value = "Thursday, 7 August, 2025 9:30 AM"

# 7. Update the totals row (now row 9).
$ws.Range("P9").Value = 38.32

Write-Host "Inserted new sale row for POWER B COMPLEX and refreshed totals/timestamp"
